# Weekly update: insert two new rows (new market date 2022-01-05) at the top
# of the "Poroto verde" / Comercializadora del Agro de Limarí block, pushing
# all the existing dated rows down by two (old row 65 -> new row 67, ...,
# old row 132 -> new row 134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before row 65 (shifts 65:132 down to 67:134).
$ws.Range("A65:A66").EntireRow.Insert()

# New row 65: Magnum, 2022-01-05
$ws.Cells.Item(65, 1).Value = 2
$ws.Cells.Item(65, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(65, 3).Value = "Coquimbo"
$ws.Cells.Item(65, 4).Value = "2022-01-05"
$ws.Cells.Item(65, 5).Value = 4
$ws.Cells.Item(65, 6).Value = 100112031
$ws.Cells.Item(65, 7).Value = "Poroto verde"
$ws.Cells.Item(65, 8).Value = "Magnum"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 600
$ws.Cells.Item(65, 11).Value = 18000
$ws.Cells.Item(65, 12).Value = 19000
$ws.Cells.Item(65, 13).Value = 18500
$ws.Cells.Item(65, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(65, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(65, 16).Value = 740
$ws.Cells.Item(65, 17).Value = 25
$ws.Cells.Item(65, 18).Value = "Hortaliza"

# New row 66: Sin especificar, 2022-01-05
$ws.Cells.Item(66, 1).Value = 2
$ws.Cells.Item(66, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = "2022-01-05"
$ws.Cells.Item(66, 5).Value = 4
$ws.Cells.Item(66, 6).Value = 100112031
$ws.Cells.Item(66, 7).Value = "Poroto verde"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 300
$ws.Cells.Item(66, 11).Value = 25000
$ws.Cells.Item(66, 12).Value = 28000
$ws.Cells.Item(66, 13).Value = 26500
$ws.Cells.Item(66, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(66, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(66, 16).Value = 1060
$ws.Cells.Item(66, 17).Value = 25
$ws.Cells.Item(66, 18).Value = "Hortaliza"
